# Insert a new data row at row 157 (pushing the existing rows 157-251 down
# to 158-252, i.e. a new daily price record was added to this weekly
# consolidation sheet). The new row carries a new observation for
# "Apio" (Primera, Pan de Azúcar origin) dated 2021-10-19 (serial 44488).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("157:157").Insert()

$ws.Cells.Item(157, 1).Value  = 3
$ws.Cells.Item(157, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(157, 3).Value  = "Coquimbo"
$ws.Cells.Item(157, 4).Value  = 44488
$ws.Cells.Item(157, 5).Value  = 5
$ws.Cells.Item(157, 6).Value  = 100112017
$ws.Cells.Item(157, 7).Value  = "Apio"
$ws.Cells.Item(157, 8).Value  = "Americana (o)"
$ws.Cells.Item(157, 9).Value  = "Primera"
$ws.Cells.Item(157, 10).Value = 120
$ws.Cells.Item(157, 11).Value = 8000
$ws.Cells.Item(157, 12).Value = 8000
$ws.Cells.Item(157, 13).Value = 8000
$ws.Cells.Item(157, 14).Value = "`$/docena de matas"
$ws.Cells.Item(157, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(157, 16).Value = 1333
$ws.Cells.Item(157, 17).Value = 6
$ws.Cells.Item(157, 18).Value = "Hortaliza"
